# Add a "prolificid" column (new column C) holding each worker's Prolific
# hash id, shift the existing columns (name, gender, matrices, race,
# mat_rank) one place to the right, and refresh the matrices/race/mat_rank
# values to match the newly-computed ranking (which also reorders Drew and
# Eli relative to one another).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the current "name" column (C) -- this
# pushes name->D, gender->E, matrices->F, race->G, mat_rank->H.
$ws.Columns("C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "prolificid"

# Prolific ids for each of the 12 workers (rows 2-13), in row order.
$prolificIds = @(
    "5e2522d6b734b47915f88275",
    "601d69a993d94008fb2b25dc",
    "60db4fde6193c50664c9c478",
    "5dd671942b033b5ec8bc97b4",
    "5ff8ad350d084e10f500e48a",
    "60b83826821417f8e484a207",
    "60b322994d0b901954690036",
    "60bf9943e4e04642d4634ecc",
    "60c2341fe95d71ee52c043f0",
    "60b091ed11ccda59e3fc7761",
    "6088fc724afd5c008db33e9d",
    "6097b95056caf5ebb2720002"
)

# Refreshed values for the shifted columns: B (offer id), D (name),
# F (matrices score) and G (race), reflecting the latest computation.
$offerIds  = @(2, 3, 22, 26, 27, 29, 33, 32, 30, 44, 49, 50)
$names     = @("Corey", "Quinterius", "Edosagbe", "Juan", "Drew", "Eli", "Brennan", "Jamarii", "Matthew", "Myles", "Masuf", "Damian")
$matrices  = @(14.40027434213815, 13.3019045588994, 8.361239919308625, 7.001638653657909, 6.38247157770629, 6.166501710522055, 5.482491394747845, 5.47701608724167, 5.11301045113269, 3.317147241751798, 1.019206744813832, 0.2975636940466398)
$races     = @("White", "Black or African American", "Black or African American", "Hispanic", "White", "White", "White", "Black or African American", "White", "Black or African American", "Asian", "Black or African American")

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $offerIds[$i]
    $ws.Cells.Item($row, 3).Value = $prolificIds[$i]
    $ws.Cells.Item($row, 4).Value = $names[$i]
    $ws.Cells.Item($row, 5).Value = "male"
    $ws.Cells.Item($row, 6).Value = $matrices[$i]
    $ws.Cells.Item($row, 7).Value = $races[$i]
    $ws.Cells.Item($row, 8).Value = $i + 1
}
